$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
  @{ Row=2; D='29.489.43'; DForce=$false; E='  +1.68%  ' },
  @{ Row=3; D='1.914.77'; DForce=$false; E='  +0.48%  ' },
  @{ Row=4; E='  -0.02%  ' },
  @{ Row=5; D='328.04'; DForce=$true; E='  -1.51%  ' },
  @{ Row=6; D='0.9999'; DForce=$true; E='  -0.02%  ' },
  @{ Row=7; D='0.4791'; DForce=$true; E='  +3.44%  ' },
  @{ Row=8; D='0.4100'; DForce=$true },
  @{ Row=9; B='OKB'; C='https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'; D='47.71'; DForce=$true; E='  -0.34%  ' },
  @{ Row=10; B='Dogecoin'; C='https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'; D='0.08027'; DForce=$true; E='  +0.13%  ' },
  @{ Row=11; B='Polygon'; C='https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'; D='1.011'; DForce=$true; E='  +0.74%  ' },
  @{ Row=12; B='Solana'; C='https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; D='22.43'; DForce=$true; E='  +3.06%  ' },
  @{ Row=13; B='WrappedEther'; C='https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; D='1.929.40'; DForce=$false; E='  +1.00%  ' },
  @{ Row=14; B='Polkadot'; C='https://coinranking.com/coin/25W7FG7om+polkadot-dot'; D='5.959'; DForce=$true; E='  +0.29%  ' },
  @{ Row=15; B='Chainlink'; C='https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; D='7.169'; DForce=$true; E='  +1.24%  ' },
  @{ Row=16; B='Litecoin'; C='https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; D='89.40'; DForce=$true; E='  +0.56%  ' },
  @{ Row=17; B='BinanceUSD'; C='https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; D='0.9996'; DForce=$true; E='  -0.10%  ' },
  @{ Row=18; B='TRON'; C='https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; D='0.06597'; DForce=$true; E='  +0.50%  ' },
  @{ Row=19; B='ShibaInu'; C='https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; D='0.00001032'; DForce=$true; E='  +0.17%  ' },
  @{ Row=20; B='Avalanche'; C='https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; D='17.77'; DForce=$true; E='  +1.56%  ' },
  @{ Row=21; B='Dai'; C='https://coinranking.com/coin/MoTuySvg7+dai-dai'; D='1.000'; DForce=$true; E='  -0.05%  ' },
  @{ Row=22; B='WrappedBTC'; C='https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; D='29.509.18'; DForce=$false; E='  +1.64%  ' },
  @{ Row=23; B='Uniswap'; C='https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; D='5.552'; DForce=$true; E='  +1.79%  ' },
  @{ Row=24; B='Cosmos'; C='https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D='11.54'; DForce=$true; E='  +2.31%  ' },
  @{ Row=25; B='Toncoin'; C='https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; D='2.205'; DForce=$true; E='  -1.51%  ' },
  @{ Row=26; B='WrappedliquidstakedEther2.0'; C='https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; D='2.145.27'; DForce=$false; E='  +0.41%  ' },
  @{ Row=27; B='Monero'; C='https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D='153.46'; DForce=$true; E='  -2.77%  ' },
  @{ Row=28; B='EthereumClassic'; C='https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D='19.81'; DForce=$true; E='  +0.37%  ' },
  @{ Row=29; B='InternetComputer(DFINITY)'; C='https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; D='5.772'; DForce=$true; E='  +6.85%  ' },
  @{ Row=30; B='LidoDAOToken'; C='https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; D='2.138'; DForce=$true; E='  +1.82%  ' },
  @{ Row=31; B='BitcoinCash'; C='https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; D='117.46'; DForce=$true; E='  -1.23%  ' },
  @{ Row=32; B='ImmutableX'; C='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D='1.067'; DForce=$true; E='  +8.88%  ' },
  @{ Row=33; B='Stellar'; C='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; D='0.09582'; DForce=$true; E='  +1.87%  ' },
  @{ Row=34; B='ARBITRUM'; C='https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; D='1.425'; DForce=$true; E='  +0.35%  ' },
  @{ Row=35; B='HuobiToken'; C='https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; D='3.572'; DForce=$true; E='  -0.53%  ' },
  @{ Row=36; B='Filecoin'; C='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D='5.394'; DForce=$true; E='  +1.53%  ' },
  @{ Row=37; B='Hedera'; C='https://coinranking.com/coin/jad286TjB+hedera-hbar'; D='0.06101'; DForce=$true; E='  +0.22%  ' },
  @{ Row=38; B='VeChain'; C='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D='0.02255'; DForce=$true; E='  +0.65%  ' },
  @{ Row=39; B='FraxShare'; C='https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; D='8.362'; DForce=$true; E='  -0.30%  ' },
  @{ Row=40; B='TrustWalletToken'; C='https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; D='1.175'; DForce=$true; E='  +0.39%  ' },
  @{ Row=41; B='TheSandbox'; C='https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; D='0.5891'; DForce=$true; E='  +1.39%  ' },
  @{ Row=42; B='Algorand'; C='https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; D='0.1840'; DForce=$true; E='  +1.05%  ' },
  @{ Row=43; B='Aptos'; C='https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; D='10.15'; DForce=$true; E='  -0.23%  ' },
  @{ Row=44; B='WEMIXTOKEN'; C='https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'; D='1.299'; DForce=$true; E='  +2.74%  ' },
  @{ Row=45; B='RenderToken'; C='https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; D='2.414'; DForce=$true; E='  +5.03%  ' },
  @{ Row=46; B='Cronos'; C='https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; D='0.07791'; DForce=$true; E='  +10.88%  ' },
  @{ Row=47; B='EnergySwap'; C='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D='12.21'; DForce=$true; E='  +0.40%  ' },
  @{ Row=48; B='Decentraland'; C='https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'; D='0.5548'; DForce=$true; E='  +0.78%  ' },
  @{ Row=49; B='NEARProtocol'; C='https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; D='1.932'; DForce=$true; E='  +1.17%  ' },
  @{ Row=50; B='Quant'; C='https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; D='113.52'; DForce=$true; E='  +1.35%  ' },
  @{ Row=51; B='Elrond'; C='https://coinranking.com/coin/omwkOTglq+elrond-egld'; D='45.38'; DForce=$true; E='  -6.32%  ' }
)

foreach ($u in $updates) {
  if ($u.ContainsKey('B')) { $ws.Cells.Item($u.Row, 2).Value = $u.B }
  if ($u.ContainsKey('C')) { $ws.Cells.Item($u.Row, 3).Value = $u.C }
  if ($u.ContainsKey('D')) {
    $cellD = $ws.Cells.Item($u.Row, 4)
    if ($u.DForce) {
      $cellD.NumberFormat = "@"
      $cellD.Value = $u.D
      $cellD.Style = "Normal"
    } else {
      $cellD.Value = $u.D
    }
  }
  if ($u.ContainsKey('E')) { $ws.Cells.Item($u.Row, 5).Value = $u.E }
}
